$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for rows with changed values.
# D-column values are forced to remain text (matching the original inlineStr
# cell type) by temporarily applying a text NumberFormat before assignment,
# then resetting the style back to Normal so no stray style index is left
# behind -- this mirrors how the source file stores these as plain strings.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.564.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.64%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.641.64"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.53%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3762"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.43%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "52.64"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.23%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3673"
$ws.Range("D9").Style = "Normal"
$ws.Range("E10").Value = "  +1.32%  "
$ws.Range("E11").Value = "  +1.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9999"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.05"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.677"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.21%  "
$ws.Range("E15").Value = "  +2.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.436"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.642.67"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "95.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06918"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.78%  "
$ws.Range("E20").Value = "  +1.73%  "
$ws.Range("E21").Value = "  +0.86%  "
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "23.558.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.090"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.90%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.419"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.67%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.58"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.360"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "135.93"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.391"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.826.46"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.71%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.824"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9799"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02849"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.80%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.45"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.07398"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.59%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2553"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.221"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.08908"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.31%  "
$ws.Range("E41").Value = "  +1.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7147"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6578"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.355"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.047"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.78%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9990"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07991"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "130.33"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.211"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.21%  "

# Rows 43 and 44 swapped coins (Aptos <-> EnergySwap) with updated price/volume.
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.40"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.92%  "
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.59"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.49%  "
